# Update NHS England deaths data (accessed 22 July 2020)
# - Revises "Daily hospital deaths" (col C) for several previously reported
#   dates; col D ("Cumulative hospital deaths") recalculates automatically
#   via its existing SUM($C$2:Cn) formula.
# - Appends newly reported daily figures for 2020-06-06 .. 2020-07-20
#   (rows 99-143), following the same A/B/C/D layout as existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed C values in existing rows (8-98) ---
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 3).Value = 4
$ws.Cells.Item(12, 3).Value = 9
$ws.Cells.Item(15, 3).Value = 22
$ws.Cells.Item(16, 3).Value = 27
$ws.Cells.Item(17, 3).Value = 40
$ws.Cells.Item(18, 3).Value = 46
$ws.Cells.Item(19, 3).Value = 65
$ws.Cells.Item(20, 3).Value = 63
$ws.Cells.Item(21, 3).Value = 105
$ws.Cells.Item(22, 3).Value = 103
$ws.Cells.Item(23, 3).Value = 149
$ws.Cells.Item(24, 3).Value = 159
$ws.Cells.Item(25, 3).Value = 204
$ws.Cells.Item(26, 3).Value = 263
$ws.Cells.Item(30, 3).Value = 437
$ws.Cells.Item(31, 3).Value = 496
$ws.Cells.Item(33, 3).Value = 644
$ws.Cells.Item(34, 3).Value = 647
$ws.Cells.Item(38, 3).Value = 726
$ws.Cells.Item(39, 3).Value = 812
$ws.Cells.Item(41, 3).Value = 790
$ws.Cells.Item(43, 3).Value = 779
$ws.Cells.Item(44, 3).Value = 717
$ws.Cells.Item(45, 3).Value = 698
$ws.Cells.Item(48, 3).Value = 639
$ws.Cells.Item(49, 3).Value = 609
$ws.Cells.Item(51, 3).Value = 522
$ws.Cells.Item(52, 3).Value = 565
$ws.Cells.Item(54, 3).Value = 501
$ws.Cells.Item(55, 3).Value = 451
$ws.Cells.Item(57, 3).Value = 385
$ws.Cells.Item(59, 3).Value = 344
$ws.Cells.Item(60, 3).Value = 341
$ws.Cells.Item(61, 3).Value = 324
$ws.Cells.Item(62, 3).Value = 312
$ws.Cells.Item(63, 3).Value = 306
$ws.Cells.Item(64, 3).Value = 268
$ws.Cells.Item(66, 3).Value = 259
$ws.Cells.Item(67, 3).Value = 251
$ws.Cells.Item(68, 3).Value = 266
$ws.Cells.Item(70, 3).Value = 213
$ws.Cells.Item(71, 3).Value = 203
$ws.Cells.Item(72, 3).Value = 196
$ws.Cells.Item(73, 3).Value = 166
$ws.Cells.Item(75, 3).Value = 162
$ws.Cells.Item(76, 3).Value = 179
$ws.Cells.Item(77, 3).Value = 171
$ws.Cells.Item(79, 3).Value = 137
$ws.Cells.Item(80, 3).Value = 158
$ws.Cells.Item(81, 3).Value = 144
$ws.Cells.Item(82, 3).Value = 153
$ws.Cells.Item(83, 3).Value = 149
$ws.Cells.Item(84, 3).Value = 121
$ws.Cells.Item(85, 3).Value = 128
$ws.Cells.Item(86, 3).Value = 116
$ws.Cells.Item(87, 3).Value = 133
$ws.Cells.Item(88, 3).Value = 138
$ws.Cells.Item(89, 3).Value = 120
$ws.Cells.Item(90, 3).Value = 124
$ws.Cells.Item(91, 3).Value = 116
$ws.Cells.Item(92, 3).Value = 92
$ws.Cells.Item(93, 3).Value = 83
$ws.Cells.Item(94, 3).Value = 94
$ws.Cells.Item(95, 3).Value = 109
$ws.Cells.Item(96, 3).Value = 110
$ws.Cells.Item(97, 3).Value = 83
$ws.Cells.Item(98, 3).Value = 86

# --- Append new rows 99-143 ---
$ws.Cells.Item(99, 1).Value = "England"
$ws.Cells.Item(99, 2).Value = 43988
$ws.Cells.Item(99, 3).Value = 83
$ws.Range("D99").Formula = "=SUM(`$C`$2:C99)"
$ws.Cells.Item(100, 1).Value = "England"
$ws.Cells.Item(100, 2).Value = 43989
$ws.Cells.Item(100, 3).Value = 80
$ws.Range("D100").Formula = "=SUM(`$C`$2:C100)"
$ws.Cells.Item(101, 1).Value = "England"
$ws.Cells.Item(101, 2).Value = 43990
$ws.Cells.Item(101, 3).Value = 73
$ws.Range("D101").Formula = "=SUM(`$C`$2:C101)"
$ws.Cells.Item(102, 1).Value = "England"
$ws.Cells.Item(102, 2).Value = 43991
$ws.Cells.Item(102, 3).Value = 67
$ws.Range("D102").Formula = "=SUM(`$C`$2:C102)"
$ws.Cells.Item(103, 1).Value = "England"
$ws.Cells.Item(103, 2).Value = 43992
$ws.Cells.Item(103, 3).Value = 77
$ws.Range("D103").Formula = "=SUM(`$C`$2:C103)"
$ws.Cells.Item(104, 1).Value = "England"
$ws.Cells.Item(104, 2).Value = 43993
$ws.Cells.Item(104, 3).Value = 49
$ws.Range("D104").Formula = "=SUM(`$C`$2:C104)"
$ws.Cells.Item(105, 1).Value = "England"
$ws.Cells.Item(105, 2).Value = 43994
$ws.Cells.Item(105, 3).Value = 52
$ws.Range("D105").Formula = "=SUM(`$C`$2:C105)"
$ws.Cells.Item(106, 1).Value = "England"
$ws.Cells.Item(106, 2).Value = 43995
$ws.Cells.Item(106, 3).Value = 43
$ws.Range("D106").Formula = "=SUM(`$C`$2:C106)"
$ws.Cells.Item(107, 1).Value = "England"
$ws.Cells.Item(107, 2).Value = 43996
$ws.Cells.Item(107, 3).Value = 58
$ws.Range("D107").Formula = "=SUM(`$C`$2:C107)"
$ws.Cells.Item(108, 1).Value = "England"
$ws.Cells.Item(108, 2).Value = 43997
$ws.Cells.Item(108, 3).Value = 56
$ws.Range("D108").Formula = "=SUM(`$C`$2:C108)"
$ws.Cells.Item(109, 1).Value = "England"
$ws.Cells.Item(109, 2).Value = 43998
$ws.Cells.Item(109, 3).Value = 60
$ws.Range("D109").Formula = "=SUM(`$C`$2:C109)"
$ws.Cells.Item(110, 1).Value = "England"
$ws.Cells.Item(110, 2).Value = 43999
$ws.Cells.Item(110, 3).Value = 50
$ws.Range("D110").Formula = "=SUM(`$C`$2:C110)"
$ws.Cells.Item(111, 1).Value = "England"
$ws.Cells.Item(111, 2).Value = 44000
$ws.Cells.Item(111, 3).Value = 49
$ws.Range("D111").Formula = "=SUM(`$C`$2:C111)"
$ws.Cells.Item(112, 1).Value = "England"
$ws.Cells.Item(112, 2).Value = 44001
$ws.Cells.Item(112, 3).Value = 42
$ws.Range("D112").Formula = "=SUM(`$C`$2:C112)"
$ws.Cells.Item(113, 1).Value = "England"
$ws.Cells.Item(113, 2).Value = 44002
$ws.Cells.Item(113, 3).Value = 45
$ws.Range("D113").Formula = "=SUM(`$C`$2:C113)"
$ws.Cells.Item(114, 1).Value = "England"
$ws.Cells.Item(114, 2).Value = 44003
$ws.Cells.Item(114, 3).Value = 36
$ws.Range("D114").Formula = "=SUM(`$C`$2:C114)"
$ws.Cells.Item(115, 1).Value = "England"
$ws.Cells.Item(115, 2).Value = 44004
$ws.Cells.Item(115, 3).Value = 42
$ws.Range("D115").Formula = "=SUM(`$C`$2:C115)"
$ws.Cells.Item(116, 1).Value = "England"
$ws.Cells.Item(116, 2).Value = 44005
$ws.Cells.Item(116, 3).Value = 51
$ws.Range("D116").Formula = "=SUM(`$C`$2:C116)"
$ws.Cells.Item(117, 1).Value = "England"
$ws.Cells.Item(117, 2).Value = 44006
$ws.Cells.Item(117, 3).Value = 54
$ws.Range("D117").Formula = "=SUM(`$C`$2:C117)"
$ws.Cells.Item(118, 1).Value = "England"
$ws.Cells.Item(118, 2).Value = 44007
$ws.Cells.Item(118, 3).Value = 46
$ws.Range("D118").Formula = "=SUM(`$C`$2:C118)"
$ws.Cells.Item(119, 1).Value = "England"
$ws.Cells.Item(119, 2).Value = 44008
$ws.Cells.Item(119, 3).Value = 36
$ws.Range("D119").Formula = "=SUM(`$C`$2:C119)"
$ws.Cells.Item(120, 1).Value = "England"
$ws.Cells.Item(120, 2).Value = 44009
$ws.Cells.Item(120, 3).Value = 30
$ws.Range("D120").Formula = "=SUM(`$C`$2:C120)"
$ws.Cells.Item(121, 1).Value = "England"
$ws.Cells.Item(121, 2).Value = 44010
$ws.Cells.Item(121, 3).Value = 37
$ws.Range("D121").Formula = "=SUM(`$C`$2:C121)"
$ws.Cells.Item(122, 1).Value = "England"
$ws.Cells.Item(122, 2).Value = 44011
$ws.Cells.Item(122, 3).Value = 29
$ws.Range("D122").Formula = "=SUM(`$C`$2:C122)"
$ws.Cells.Item(123, 1).Value = "England"
$ws.Cells.Item(123, 2).Value = 44012
$ws.Cells.Item(123, 3).Value = 29
$ws.Range("D123").Formula = "=SUM(`$C`$2:C123)"
$ws.Cells.Item(124, 1).Value = "England"
$ws.Cells.Item(124, 2).Value = 44013
$ws.Cells.Item(124, 3).Value = 17
$ws.Range("D124").Formula = "=SUM(`$C`$2:C124)"
$ws.Cells.Item(125, 1).Value = "England"
$ws.Cells.Item(125, 2).Value = 44014
$ws.Cells.Item(125, 3).Value = 33
$ws.Range("D125").Formula = "=SUM(`$C`$2:C125)"
$ws.Cells.Item(126, 1).Value = "England"
$ws.Cells.Item(126, 2).Value = 44015
$ws.Cells.Item(126, 3).Value = 18
$ws.Range("D126").Formula = "=SUM(`$C`$2:C126)"
$ws.Cells.Item(127, 1).Value = "England"
$ws.Cells.Item(127, 2).Value = 44016
$ws.Cells.Item(127, 3).Value = 21
$ws.Range("D127").Formula = "=SUM(`$C`$2:C127)"
$ws.Cells.Item(128, 1).Value = "England"
$ws.Cells.Item(128, 2).Value = 44017
$ws.Cells.Item(128, 3).Value = 23
$ws.Range("D128").Formula = "=SUM(`$C`$2:C128)"
$ws.Cells.Item(129, 1).Value = "England"
$ws.Cells.Item(129, 2).Value = 44018
$ws.Cells.Item(129, 3).Value = 23
$ws.Range("D129").Formula = "=SUM(`$C`$2:C129)"
$ws.Cells.Item(130, 1).Value = "England"
$ws.Cells.Item(130, 2).Value = 44019
$ws.Cells.Item(130, 3).Value = 22
$ws.Range("D130").Formula = "=SUM(`$C`$2:C130)"
$ws.Cells.Item(131, 1).Value = "England"
$ws.Cells.Item(131, 2).Value = 44020
$ws.Cells.Item(131, 3).Value = 21
$ws.Range("D131").Formula = "=SUM(`$C`$2:C131)"
$ws.Cells.Item(132, 1).Value = "England"
$ws.Cells.Item(132, 2).Value = 44021
$ws.Cells.Item(132, 3).Value = 37
$ws.Range("D132").Formula = "=SUM(`$C`$2:C132)"
$ws.Cells.Item(133, 1).Value = "England"
$ws.Cells.Item(133, 2).Value = 44022
$ws.Cells.Item(133, 3).Value = 16
$ws.Range("D133").Formula = "=SUM(`$C`$2:C133)"
$ws.Cells.Item(134, 1).Value = "England"
$ws.Cells.Item(134, 2).Value = 44023
$ws.Cells.Item(134, 3).Value = 9
$ws.Range("D134").Formula = "=SUM(`$C`$2:C134)"
$ws.Cells.Item(135, 1).Value = "England"
$ws.Cells.Item(135, 2).Value = 44024
$ws.Cells.Item(135, 3).Value = 14
$ws.Range("D135").Formula = "=SUM(`$C`$2:C135)"
$ws.Cells.Item(136, 1).Value = "England"
$ws.Cells.Item(136, 2).Value = 44025
$ws.Cells.Item(136, 3).Value = 21
$ws.Range("D136").Formula = "=SUM(`$C`$2:C136)"
$ws.Cells.Item(137, 1).Value = "England"
$ws.Cells.Item(137, 2).Value = 44026
$ws.Cells.Item(137, 3).Value = 10
$ws.Range("D137").Formula = "=SUM(`$C`$2:C137)"
$ws.Cells.Item(138, 1).Value = "England"
$ws.Cells.Item(138, 2).Value = 44027
$ws.Cells.Item(138, 3).Value = 17
$ws.Range("D138").Formula = "=SUM(`$C`$2:C138)"
$ws.Cells.Item(139, 1).Value = "England"
$ws.Cells.Item(139, 2).Value = 44028
$ws.Cells.Item(139, 3).Value = 8
$ws.Range("D139").Formula = "=SUM(`$C`$2:C139)"
$ws.Cells.Item(140, 1).Value = "England"
$ws.Cells.Item(140, 2).Value = 44029
$ws.Cells.Item(140, 3).Value = 6
$ws.Range("D140").Formula = "=SUM(`$C`$2:C140)"
$ws.Cells.Item(141, 1).Value = "England"
$ws.Cells.Item(141, 2).Value = 44030
$ws.Cells.Item(141, 3).Value = 10
$ws.Range("D141").Formula = "=SUM(`$C`$2:C141)"
$ws.Cells.Item(142, 1).Value = "England"
$ws.Cells.Item(142, 2).Value = 44031
$ws.Cells.Item(142, 3).Value = 4
$ws.Range("D142").Formula = "=SUM(`$C`$2:C142)"
$ws.Cells.Item(143, 1).Value = "England"
$ws.Cells.Item(143, 2).Value = 44032
$ws.Cells.Item(143, 3).Value = 0
$ws.Range("D143").Formula = "=SUM(`$C`$2:C143)"

# --- Update selection ---
$ws.Range("G29").Select()
